$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the old extent first (rows 1-5, cols A-B) so stale cells (A4, B4, B5) go away.
$ws.Range("A1:B5").Clear()

# Row 1 - headers
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "Invalid UserName"
$ws.Range("D1").Value = "Invalid Password"
$ws.Range("E1").Value = "Appointment Date"

# Row 2 - data
$ws.Range("A2").Value = "John Doe"
$ws.Range("B2").Value = "ThisIsNotAPassword"
$ws.Range("C2").Value = "John DO"
$ws.Range("D2").Value = "ThisIsNotPass"
$ws.Range("E2").Value = (Get-Date -Year 2022 -Month 5 -Day 3)
$ws.Range("E2").NumberFormat = "m/d/yyyy"

# Row 3 - data
$ws.Range("A3").Value = "John Test1"
$ws.Range("B3").Value = "Password"

# Column widths
$ws.Range("A1").EntireColumn.ColumnWidth = 20.42578125
$ws.Range("B1").EntireColumn.ColumnWidth = 24.5703125
$ws.Range("C1").EntireColumn.ColumnWidth = 16.85546875
$ws.Range("D1").EntireColumn.ColumnWidth = 21
$ws.Range("E1").EntireColumn.ColumnWidth = 20.85546875

# Hyperlink at A5 (kept from original template, cell itself has no value)
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:John@", "", "", "John@")

# Selection
$ws.Range("C3").Select()
